# NCAAF Futures - "week 2 excel file" upload
# Mirrors the author's edit:
#   1. Rename the (only) worksheet from "Sheet2" to "Sheet1"
#   2. Clean up a handful of floating point rounding artifacts in column C
#      (e.g. 1.0004999999999999 / 0.99950000000000006 -> 1)
#   3. Fix a sort-order swap between two pairs of teams that tied on odds
#   4. Append the Week 2 odds table (50 new rows, Week = 2) below the
#      existing Week 1 table
#   5. Leave the selection/scroll where the author left it after pasting
#      the new data (near the bottom of the new block)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename sheet ------------------------------------------------------
$ws.Name = "Sheet1"

# --- 2. Round the odds-ratio noise down to a clean 1 ----------------------
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("C20").Value = 1

# --- 3. Swap the two tied-odds pairs back into the correct order ----------
$ws.Range("B35").Value = "Colorado"
$ws.Range("B36").Value = "TCU"

$ws.Range("B42").Value = "South Carolina"
$ws.Range("B43").Value = "Iowa State"

# --- 4. Append the Week 2 odds table (rows 52-101) -------------------------
$week2 = @(
    @{Row=52; Odds=280; Team='Georgia'},
    @{Row=53; Odds=400; Team='Ohio State'},
    @{Row=54; Odds=700; Team='Oregon'},
    @{Row=55; Odds=700; Team='Texas'},
    @{Row=56; Odds=1200; Team='Alabama'},
    @{Row=57; Odds=1400; Team='Ole Miss'},
    @{Row=58; Odds=1600; Team='Notre Dame'},
    @{Row=59; Odds=1800; Team='Penn State'},
    @{Row=60; Odds=2500; Team='Miami'},
    @{Row=61; Odds=3500; Team='Michigan'},
    @{Row=62; Odds=3000; Team='Tennessee'},
    @{Row=63; Odds=3500; Team='LSU'},
    @{Row=64; Odds=5000; Team='Utah'},
    @{Row=65; Odds=4000; Team='Missouri'},
    @{Row=66; Odds=6000; Team='Clemson'},
    @{Row=67; Odds=5000; Team='USC'},
    @{Row=68; Odds=6000; Team='Texas A&M'},
    @{Row=69; Odds=6000; Team='Oklahoma'},
    @{Row=70; Odds=6600; Team='Kansas State'},
    @{Row=71; Odds=12500; Team='Auburn'},
    @{Row=72; Odds=20000; Team='Colorado'},
    @{Row=73; Odds=12500; Team='NC State'},
    @{Row=74; Odds=20000; Team='Louisville'},
    @{Row=75; Odds=25000; Team='Kansas'},
    @{Row=76; Odds=20000; Team='Arizona'},
    @{Row=77; Odds=20000; Team='Nebraska'},
    @{Row=78; Odds=12500; Team='Oklahoma State'},
    @{Row=79; Odds=20000; Team='Iowa'},
    @{Row=80; Odds=25000; Team='Wisconsin'},
    @{Row=81; Odds=40000; Team='Boise State'},
    @{Row=82; Odds=30000; Team='SMU'},
    @{Row=83; Odds=30000; Team='Washington'},
    @{Row=84; Odds=35000; Team='TCU'},
    @{Row=85; Odds=30000; Team='Kentucky'},
    @{Row=86; Odds=25000; Team='Texas Tech'},
    @{Row=87; Odds=40000; Team='Florida State'},
    @{Row=88; Odds=50000; Team='Oregon State'},
    @{Row=89; Odds=30000; Team='Iowa State'},
    @{Row=90; Odds=30000; Team='UCF'},
    @{Row=91; Odds=40000; Team='Virginia Tech'},
    @{Row=92; Odds=30000; Team='Florida'},
    @{Row=93; Odds=60000; Team='Tulane'},
    @{Row=94; Odds=40000; Team='North Carolina'},
    @{Row=95; Odds=100000; Team='UCLA'},
    @{Row=96; Odds=40000; Team='Arkansas'},
    @{Row=97; Odds=50000; Team='Memphis'},
    @{Row=98; Odds=35000; Team='Georgia Tech'},
    @{Row=99; Odds=60000; Team='Liberty'},
    @{Row=100; Odds=40000; Team='South Carolina'},
    @{Row=101; Odds=40000; Team='Maryland'}
)

foreach ($row in $week2) {
    $ws.Cells.Item($row.Row, 1).Value = $row.Odds
    $ws.Cells.Item($row.Row, 2).Value = $row.Team
    $ws.Cells.Item($row.Row, 3).Value = 2
}

# --- 5. Match the author's final selection/scroll position ----------------
$excel.ActiveWindow.ScrollRow = 73
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("A52:B101").Select()
